$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data scraped on Wed Nov  8 20:44:15 UTC 2023

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.534.65"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.893.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.691"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "43.02"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "56.97"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +9.39%  "
$ws.Range("E10").Value = "  +0.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0752"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.17%  "
$ws.Range("E12").Value = "  +1.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.56"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +10.47%  "
$ws.Range("E14").Value = "  +9.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.172.31"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.03"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.69%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.910.86"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "35.534.68"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0830"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "246.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.97"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.30%  "
$ws.Range("E23").Value = "  +3.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.58"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.65"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.35"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.45%  "
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0607"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.25"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.85"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +17.28%  "
$ws.Range("B35").Value = "BinanceUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.47"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -17.36%  "
$ws.Range("E37").Value = "  -0.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.95"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.69%  "
$ws.Range("E39").Value = "  +7.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0228"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.97"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.73%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.05"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +17.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.317.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.53%  "
$ws.Range("E46").Value = "  -1.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0811"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.45%  "
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("E49").Value = "  -0.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.40"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "42.55"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.49%  "
